$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 218, shifting all existing rows (218-263)
# down to (220-265). This mirrors the diff, which grows the used range
# from A1:R263 to A1:R265 by adding two new weekly price records (date
# 2023-06-16 / serial 45093) for "Primera" and "Segunda" quality grades.
$ws.Rows.Item(218).Insert()
$ws.Rows.Item(218).Insert()

# New row 218: Cilantro, Primera, date 2023-06-16
$ws.Range("A218").Value = 7
$ws.Range("B218").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C218").Value = "Ñuble"
$ws.Range("D218").Value = 45093
$ws.Range("E218").Value = 16
$ws.Range("F218").Value = 100112040
$ws.Range("G218").Value = "Cilantro"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 200
$ws.Range("K218").Value = 1000
$ws.Range("L218").Value = 1200
$ws.Range("M218").Value = 1100
$ws.Range("N218").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O218").Value = "Provincia de Diguillín"
$ws.Range("P218").Value = 1100
$ws.Range("Q218").Value = 1
$ws.Range("R218").Value = "Hortaliza"

# New row 219: Cilantro, Segunda, date 2023-06-16
$ws.Range("A219").Value = 7
$ws.Range("B219").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C219").Value = "Ñuble"
$ws.Range("D219").Value = 45093
$ws.Range("E219").Value = 16
$ws.Range("F219").Value = 100112040
$ws.Range("G219").Value = "Cilantro"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Segunda"
$ws.Range("J219").Value = 150
$ws.Range("K219").Value = 800
$ws.Range("L219").Value = 800
$ws.Range("M219").Value = 800
$ws.Range("N219").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O219").Value = "Provincia de Diguillín"
$ws.Range("P219").Value = 800
$ws.Range("Q219").Value = 1
$ws.Range("R219").Value = "Hortaliza"
